# Commit: Change "AWS" to "Azure" in PennDOT section
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
# MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace
$find.Execute("AWS", $true, $true, $false, $false, $false,
              $true, 1, $false, "Azure", 2)
